$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data updates (report date changed from 2018-12-31 to 2017-12-31,
# along with all the corresponding financial figures).
$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 826404898.78
$ws.Range("P2").Value = 114345838.01
$ws.Range("Q2").Value = 117257669.33
$ws.Range("R2").Value = 88.38266733179999
$ws.Range("S2").Value = 246457767.07
$ws.Range("T2").Value = 10.4417865459
$ws.Range("U2").Value = 100892673.92
$ws.Range("V2").Value = 135.1735722232
$ws.Range("W2").Value = 372712913.3
$ws.Range("X2").Value = 208681726.74
$ws.Range("Y2").Value = 217.8939041073

# These two were previously empty (inlineStr placeholder) and now get numeric values.
$ws.Range("Z2").Value = 319519.13
$ws.Range("AA2").Value = 276.2071966306

$ws.Range("AB2").Value = 453691985.48
$ws.Range("AC2").Value = 8.947100836900001
$ws.Range("AD2").Value = 41.6620907491
$ws.Range("AE2").Value = 123.2747094657
$ws.Range("AF2").Value = 215.5822936784
$ws.Range("AG2").Value = 45.1005208041
